# Insert a new weekly record row at row 124 (shifting existing rows 124-193 down to 125-194)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("124:124").Insert()

$ws.Range("A124").Value2 = 1
$ws.Range("B124").Value2 = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C124").Value2 = 'Arica y Parinacota'
$ws.Range("D124").Value2 = 44603
$ws.Range("E124").Value2 = 15
$ws.Range("F124").Value2 = 'Fruta'
$ws.Range("G124").Value2 = 100108
$ws.Range("H124").Value2 = 'Tropicales y subtropicales'
$ws.Range("I124").Value2 = 100108006
$ws.Range("J124").Value2 = 'Plátano'
$ws.Range("K124").Value2 = 'Sin especificar'
$ws.Range("L124").Value2 = 'Pintón'
$ws.Range("M124").Value2 = 160
$ws.Range("N124").Value2 = 18000
$ws.Range("O124").Value2 = 19000
$ws.Range("P124").Value2 = 18500
$ws.Range("Q124").Value2 = '$/caja 20 kilos'
$ws.Range("R124").Value2 = 'Ecuador'
$ws.Range("S124").Value2 = 925
$ws.Range("T124").Value2 = 20
